$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1681.25
$ws.Range("I70").Value = 1500
$ws.Range("J70").Value = 1741.6666
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 5224.9998
$ws.Range("M70").Value = -4230
$ws.Range("N70").Value = -5764.9998

$ws.Range("H73").Value = 1681.25
$ws.Range("I73").Value = 1500
$ws.Range("J73").Value = 1741.6666
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 5224.9998
$ws.Range("M73").Value = -3564
$ws.Range("N73").Value = -7096.9998

$ws.Range("H80").Value = 486.77777
$ws.Range("J80").Value = 496.2
$ws.Range("L80").Value = 1488.6
$ws.Range("N80").Value = -3484.6

$ws.Range("H83").Value = 486.77777
$ws.Range("J83").Value = 496.2
$ws.Range("L83").Value = 4465.8
$ws.Range("N83").Value = -14449.8

$ws.Range("H88").Value = 3950.5
$ws.Range("J88").Value = 6000
$ws.Range("L88").Value = 6000
$ws.Range("N88").Value = -6812

$ws.Range("H91").Value = 3950.5
$ws.Range("J91").Value = 6000
$ws.Range("L91").Value = 6000
$ws.Range("N91").Value = -8808

$ws.Range("H92").Value = 930.7273
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H94").Value = 4541.125
$ws.Range("I94").Value = 4541.125
$ws.Range("K94").Value = 4541.125
$ws.Range("M94").Value = -4090.125

$ws.Range("H100").Value = 1365.4286
$ws.Range("I100").Value = 967.2
$ws.Range("K100").Value = 967.2
$ws.Range("M100").Value = -426.2

$ws.Range("H135").Value = 2664.2856
$ws.Range("J135").Value = 3050
$ws.Range("L135").Value = 27450
$ws.Range("N135").Value = -32520

$ws.Range("H138").Value = 2268.7693
$ws.Range("I138").Value = 998
$ws.Range("J138").Value = 2650
$ws.Range("K138").Value = 2994
$ws.Range("L138").Value = 7950
$ws.Range("M138").Value = 2146
$ws.Range("N138").Value = -18230

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7104.2856
$ws.Range("I32").Value = 2678.182
$ws.Range("K32").Value = 2678.182
$ws.Range("M32").Value = -2391.182

$ws.Range("H74").Value = 2557.25

$ws.Range("H77").Value = 2557.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3386.625
$ws.Range("I94").Value = 699.6667
$ws.Range("K94").Value = 699.6667
$ws.Range("M94").Value = -248.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6475
$ws.Range("J16").Value = 9950
$ws.Range("L16").Value = 9950
$ws.Range("N16").Value = -10524

$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()

$ws.Range("H88").Value = 12936
$ws.Range("J88").Value = 12936
$ws.Range("L88").Value = 12936
$ws.Range("N88").Value = -13748

$ws.Range("H91").Value = 12936
$ws.Range("J91").Value = 12936
$ws.Range("L91").Value = 12936
$ws.Range("N91").Value = -15744

$ws.Range("H94").Value = 3830
$ws.Range("I94").Value = 4537.5
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 4537.5
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -4086.5
$ws.Range("N94").Value = -1902

$ws.Range("H99").Value = 6048.75
$ws.Range("I99").Value = 2897.5
$ws.Range("K99").Value = 2897.5
$ws.Range("M99").Value = -1399.5

$ws.Range("H113").Value = 6475
$ws.Range("J113").Value = 9950
$ws.Range("L113").Value = 9950
$ws.Range("N113").Value = -14290

$ws.Range("H126").Value = 6048.75
$ws.Range("I126").Value = 2897.5
$ws.Range("K126").Value = 8692.5
$ws.Range("M126").Value = -6222.5

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 80.40000000000001
$ws.Range("I2").Value = 134.66667
$ws.Range("K2").Value = 808.0000200000001
$ws.Range("M2").Value = -695.0000200000001

$ws.Range("H17").Value = 694
$ws.Range("I17").Value = 79.5
$ws.Range("J17").Value = 1308.5
$ws.Range("K17").Value = 238.5
$ws.Range("L17").Value = 3925.5
$ws.Range("M17").Value = -69.5
$ws.Range("N17").Value = -4263.5

$ws.Range("H39").Value = 1300
$ws.Range("J39").Value = 2500
$ws.Range("L39").Value = 7500
$ws.Range("N39").Value = -8088

$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H55").Value = 1530
$ws.Range("J55").Value = 1845
$ws.Range("L55").Value = 5535
$ws.Range("N55").Value = -5889

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H97").Value = 272
$ws.Range("J97").Value = 125
$ws.Range("L97").Value = 375
$ws.Range("N97").Value = -1367

$ws.Range("H139").Value = 1989.1428
$ws.Range("I139").Value = 1989.1428
$ws.Range("K139").Value = 5967.428400000001
$ws.Range("M139").Value = -827.4284000000007

$ws.Range("H140").Value = 535.2222
$ws.Range("I140").Value = 535.2222
$ws.Range("K140").Value = 1605.6666
$ws.Range("M140").Value = 3574.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws.Range("H93").Value = 2000
$ws.Range("I93").Value = 2000
$ws.Range("K93").Value = 2000
$ws.Range("M93").Value = -752

$ws.Range("H94").Value = 61473
$ws.Range("J94").Value = 61473
$ws.Range("L94").Value = 61473
$ws.Range("N94").Value = -62825

$ws.Range("H101").Value = 17499
$ws.Range("J101").Value = 17499
$ws.Range("L101").Value = 17499
$ws.Range("N101").Value = -23989

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 1500
$ws.Range("J21").Value = 1500
$ws.Range("L21").Value = 1500
$ws.Range("N21").Value = -1970

$ws.Range("H35").Value = 1500
$ws.Range("J35").Value = 1500
$ws.Range("L35").Value = 1500
$ws.Range("N35").Value = -2080

$ws.Range("H81").Value = 999
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()

$ws.Range("H84").Value = 999
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()

$ws.Range("H100").Value = 575
$ws.Range("I100").Value = 575
$ws.Range("K100").Value = 1150
$ws.Range("M100").Value = -609

$ws.Range("H113").Value = 513.3333
$ws.Range("I113").Value = 395
$ws.Range("K113").Value = 1185
$ws.Range("M113").Value = 985

$ws.Range("H131").Value = 97999.5
$ws.Range("J131").Value = 97999.5
$ws.Range("L131").Value = 97999.5
$ws.Range("N131").Value = -108079.5

$ws.Range("H132").Value = 1997.6666
